$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

$ws.Range('D2').Value = '66.887.46'
$ws.Range('E2').Value = '  -0.87%  '
$ws.Range('D3').Value = '3.463.46'
$ws.Range('E3').Value = '  -1.94%  '
$ws.Range('E4').Value = '  -0.02%  '
Set-TextValue 'D5' '592.40'
$ws.Range('E5').Value = '  -0.80%  '
Set-TextValue 'D6' '174.98'
$ws.Range('E6').Value = '  +0.77%  '
Set-TextValue 'D7' '0.999'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  -1.74%  '
Set-TextValue 'D9' '0.129'
$ws.Range('E9').Value = '  -4.70%  '
Set-TextValue 'D10' '7.07'
$ws.Range('E10').Value = '  -3.42%  '
$ws.Range('E11').Value = '  -3.04%  '
$ws.Range('D12').Value = '4.060.81'
$ws.Range('E12').Value = '  -2.03%  '
Set-TextValue 'D13' '30.69'
$ws.Range('E13').Value = '  +6.04%  '
$ws.Range('E14').Value = '  -0.17%  '
$ws.Range('D15').Value = '66.886.32'
$ws.Range('E15').Value = '  -0.71%  '
$ws.Range('E16').Value = '  -4.44%  '
$ws.Range('D17').Value = '3.465.90'
$ws.Range('E17').Value = '  -1.76%  '
$ws.Range('E18').Value = '  -2.60%  '
Set-TextValue 'D19' '14.27'
$ws.Range('E19').Value = '  +0.42%  '
Set-TextValue 'D20' '385.33'
$ws.Range('E20').Value = '  -3.14%  '
Set-TextValue 'D21' '7.82'
$ws.Range('E21').Value = '  -2.33%  '
Set-TextValue 'D22' '72.41'
$ws.Range('E22').Value = '  -1.60%  '
Set-TextValue 'D23' '0.995'
$ws.Range('E23').Value = '  -0.39%  '
$ws.Range('E24').Value = '  -0.23%  '
Set-TextValue 'D25' '0.531'
$ws.Range('E25').Value = '  -1.84%  '
Set-TextValue 'D26' '0.0000121'
$ws.Range('E26').Value = '  -2.03%  '
Set-TextValue 'D27' '10.27'
$ws.Range('E27').Value = '  -0.27%  '
$ws.Range('E28').Value = '  -2.51%  '
Set-TextValue 'D29' '0.996'
$ws.Range('E29').Value = '  -0.19%  '
Set-TextValue 'D30' '6.07'
$ws.Range('E30').Value = '  -3.60%  '
Set-TextValue 'D31' '1.41'
$ws.Range('E31').Value = '  -4.12%  '
$ws.Range('E32').Value = '  -2.71%  '
Set-TextValue 'D33' '23.33'
$ws.Range('E33').Value = '  -3.33%  '
Set-TextValue 'D34' '7.22'
$ws.Range('E34').Value = '  -2.44%  '
$ws.Range('E35').Value = '  -2.15%  '
Set-TextValue 'D36' '163.14'
$ws.Range('E36').Value = '  -0.46%  '
$ws.Range('E37').Value = '  -3.30%  '
$ws.Range('E38').Value = '  -0.97%  '
Set-TextValue 'D39' '7.02'
$ws.Range('E39').Value = '  +0.62%  '
$ws.Range('E40').Value = '  -0.80%  '
$ws.Range('E41').Value = '  -3.41%  '
Set-TextValue 'D42' '26.18'
$ws.Range('E42').Value = '  -1.70%  '
$ws.Range('D43').Value = '2.777.89'
$ws.Range('E43').Value = '  -1.17%  '
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue 'D44' '2.57'
$ws.Range('E44').Value = '  -2.66%  '
$ws.Range('B45').Value = 'Hedera'
$ws.Range('C45').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D45' '0.0717'
$ws.Range('E45').Value = '  -4.41%  '
Set-TextValue 'D46' '42.11'
$ws.Range('E46').Value = '  -2.01%  '
$ws.Range('E47').Value = '  -4.90%  '
Set-TextValue 'D48' '336.67'
$ws.Range('E48').Value = '  -1.69%  '
$ws.Range('E49').Value = '  -3.68%  '
Set-TextValue 'D50' '33.06'
$ws.Range('E50').Value = '  -2.52%  '
$ws.Range('E51').Value = '  -3.36%  '
